$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "D4" = "Florian"
    "D6" = "Mythron"
    "D28" = "Aquari"
    "D30" = "Sentire"
    "D35" = "Sentire"
    "D47" = "Emotivor"
    "D52" = "Emotivor"
    "D53" = "Zorblax"
    "D61" = "Nexoon"
    "D62" = "Nexoon"
    "D85" = "Aquari"
    "D86" = "Nexoon"
    "D93" = "Faerix"
    "D95" = "Zorblax"
    "D98" = "Mythron"
    "D102" = "Aquari"
    "D107" = "Quixnar"
    "D112" = "Emotivor"
    "D124" = "Nexoon"
    "D142" = "Cybex"
    "D144" = "Aquari"
    "D147" = "Florian"
    "D152" = "Florian"
    "D153" = "Nexoon"
    "D159" = "Sentire"
    "D169" = "Cybex"
    "D173" = "Sentire"
    "D177" = "Sentire"
    "D183" = "Quixnar"
    "D189" = "Nexoon"
    "D213" = "Mythron"
    "D215" = "Nexoon"
    "D241" = "Quixnar"
    "D245" = "Emotivor"
    "D249" = "Sentire"
    "D252" = "Zorblax"
    "D262" = "Faerix"
    "D265" = "Mythron"
    "D271" = "Zorblax"
    "D277" = "Quixnar"
    "D278" = "Zorblax"
    "D282" = "Sentire"
    "D296" = "Faerix"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
